# Adds the new "f2031f0a-2b1a-493f-893f-f63f1f1858ce" file's handoff/handback
# report as row 9 on all three sheets (Overview, zh-cn, de-de), matching the
# "Generate Report for Handoff" run that produced the new data rows + widened
# tables (A1:P8 -> A1:P9 / A1:G8 -> A1:G9).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Cells.Item(9, 1).Value = "f2031f0a-2b1a-493f-893f-f63f1f1858ce.md"

$ws.Cells.Item(9, 2).Value = "e2e\f2031f0a-2b1a-493f-893f-f63f1f1858ce.md"
$ws.Cells.Item(9, 2).Style = "HyperLink"

$ws.Cells.Item(9, 3).Value = ".md"
$ws.Cells.Item(9, 5).Value = "Ready for handoff"
$ws.Cells.Item(9, 6).Value = "Ready for handoff"

$ws.Cells.Item(9, 7).Value = "2016-09-06 17:18:55"
$ws.Cells.Item(9, 7).NumberFormat = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Cells.Item(9, 1).Value = "f2031f0a-2b1a-493f-893f-f63f1f1858ce.md"
$ws.Cells.Item(9, 1).Style = "HyperLink"

$ws.Cells.Item(9, 2).Value = ".md"
$ws.Cells.Item(9, 3).Value = "Ready for handoff"
$ws.Cells.Item(9, 4).Value = "e2e"
$ws.Cells.Item(9, 5).Value = "ht"
$ws.Cells.Item(9, 6).Value = "False"

$ws.Cells.Item(9, 7).Value = "f2031f0a-2b1a-493f-893f-f63f1f1858ce.80888a3f371db147e2f85912f9532ae6b5ca5e8a.zh-cn.xlf"

$ws.Cells.Item(9, 8).Value = "2016-09-06 17:18:49"
$ws.Cells.Item(9, 8).NumberFormat = "yyyy-mm-dd HH:mm:ss"

$ws.Cells.Item(9, 11).Value = "0001-01-01 00:00:00"
$ws.Cells.Item(9, 11).NumberFormat = "yyyy-mm-dd HH:mm:ss"

$ws.Cells.Item(9, 13).Value = "True"
$ws.Cells.Item(9, 15).Value = "False"

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Cells.Item(9, 1).Value = "f2031f0a-2b1a-493f-893f-f63f1f1858ce.md"
$ws.Cells.Item(9, 1).Style = "HyperLink"

$ws.Cells.Item(9, 2).Value = ".md"
$ws.Cells.Item(9, 3).Value = "Ready for handoff"
$ws.Cells.Item(9, 4).Value = "e2e"
$ws.Cells.Item(9, 5).Value = "ht"
$ws.Cells.Item(9, 6).Value = "False"

$ws.Cells.Item(9, 7).Value = "f2031f0a-2b1a-493f-893f-f63f1f1858ce.80888a3f371db147e2f85912f9532ae6b5ca5e8a.de-de.xlf"

$ws.Cells.Item(9, 8).Value = "2016-09-06 17:18:55"
$ws.Cells.Item(9, 8).NumberFormat = "yyyy-mm-dd HH:mm:ss"

$ws.Cells.Item(9, 11).Value = "0001-01-01 00:00:00"
$ws.Cells.Item(9, 11).NumberFormat = "yyyy-mm-dd HH:mm:ss"

$ws.Cells.Item(9, 13).Value = "True"
$ws.Cells.Item(9, 15).Value = "False"
